# Refresh the crypto price/volume table with the latest scraped values.
# Note: Price cells in column D are plain text (e.g. "72.224.24", "1.00"),
# not numbers. Where the new value would otherwise be auto-parsed by Excel
# as a number (losing formatting such as trailing zeros), it is written
# with a leading apostrophe to force literal text entry, matching the
# original text storage of that cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.224.24"
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("D3").Value = "2.643.80"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'602.99"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").Value = "'180.09"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("E9").Value = "  +4.77%  "
$ws.Range("D10").Value = "2.643.27"
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("D12").Value = "'0.359"
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("E14").Value = "  +3.67%  "
$ws.Range("D15").Value = "3.111.87"
$ws.Range("D16").Value = "72.175.98"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("D17").Value = "'26.64"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "2.650.29"
$ws.Range("E19").Value = "  +4.83%  "
$ws.Range("D20").Value = "'379.78"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").Value = "'7.94"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'2.08"
$ws.Range("E23").Value = "  +10.65%  "
$ws.Range("D24").Value = "'73.40"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").Value = "'4.40"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  +4.23%  "
$ws.Range("D28").Value = "2.782.75"
$ws.Range("E28").Value = "  +1.55%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("D31").Value = "'524.67"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").Value = "'8.15"
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "'164.29"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "'19.37"
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("D38").Value = "'0.113"
$ws.Range("E38").Value = "  -5.73%  "
$ws.Range("E39").Value = "  +2.18%  "
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("D41").Value = "'1.85"
$ws.Range("E41").Value = "  +1.16%  "
$ws.Range("E42").Value = "  +3.87%  "
$ws.Range("D43").Value = "'5.08"
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "'0.334"
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("D46").Value = "'39.38"
$ws.Range("E46").Value = "  -2.78%  "
$ws.Range("D47").Value = "'151.57"
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("D48").Value = "'3.72"
$ws.Range("E48").Value = "  +1.54%  "
$ws.Range("E49").Value = "  +2.65%  "
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("E51").Value = "  -2.34%  "
